$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the existing data (header row 1, data rows 2..34).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row

# --- Header row: new columns I ("I0") and J ("IF") -------------------------
# Copy H1's formatting (bold, bordered, centered header style) onto I1:J1
# before writing the header text, so the new header cells share H1's style
# (same cellXf as the rest of row 1) instead of creating a new one.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows: I = 1 (constant), J = same value as column H ---------------
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value2
}

$ws.Range("A1").Select() | Out-Null
